$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "P_1092"
$ws.Range("C2").Value = 3521.316740009102
$ws.Range("D2").Value = 156.3399505236742
$ws.Range("E2").Value = 91.07593037516521
$ws.Range("F2").Value = 7090.247189282279

$ws.Range("B3").Value = "P_1307"
$ws.Range("C3").Value = 1834.684707625472
$ws.Range("D3").Value = 87.30917999999996
$ws.Range("E3").Value = 53.283533436125
$ws.Range("F3").Value = 3959.599999999998

$ws.Range("B4").Value = "P_1111"
$ws.Range("C4").Value = 1873.177264422237
$ws.Range("D4").Value = 78.49733850000007
$ws.Range("E4").Value = 48.64042827868546
$ws.Range("F4").Value = 3559.970000000003

$ws.Range("B5").Value = "P_1393"
$ws.Range("C5").Value = 2220.28185576255
$ws.Range("D5").Value = 122.03978091017
$ws.Range("E5").Value = 52.60679542982393
$ws.Range("F5").Value = 5534.68394150431

$ws.Range("B6").Value = "P_1304"
$ws.Range("C6").Value = 1340.77721583532
$ws.Range("D6").Value = 61.56360000000002
$ws.Range("E6").Value = 36.57289490219227
$ws.Range("F6").Value = 2792

$ws.Range("B7").Value = "P_1279"
$ws.Range("C7").Value = 3736.991651557322
$ws.Range("D7").Value = 168.40908
$ws.Range("E7").Value = 96.05876624089991
$ws.Range("F7").Value = 7637.599999999999

$ws.Range("B8").Value = "P_1371"
$ws.Range("C8").Value = 1471.011631246981
$ws.Range("D8").Value = 56.64827908777654
$ws.Range("E8").Value = 33.13819390409638
$ws.Range("F8").Value = 2569.0829518266

$ws.Range("B9").Value = "P_1419"
$ws.Range("C9").Value = 1583.030761262681
$ws.Range("D9").Value = 67.53164665920092
$ws.Range("E9").Value = 38.6483330356762
$ws.Range("F9").Value = 3062.659712435416

$ws.Range("B10").Value = "P_1421"
$ws.Range("C10").Value = 1500.40559338421
$ws.Range("D10").Value = 65.61978126510103
$ws.Range("E10").Value = 36.30818143506726
$ws.Range("F10").Value = 2975.953798870795

$ws.Range("B11").Value = "P_1100"
$ws.Range("C11").Value = 3455.249038955813
$ws.Range("D11").Value = 149.7956827499998
$ws.Range("E11").Value = 89.23199304412083
$ws.Range("F11").Value = 6793.454999999992

$ws.Range("B12").Value = "P_1141"
$ws.Range("C12").Value = 4647.068705195382
$ws.Range("D12").Value = 178.5829499999998
$ws.Range("E12").Value = 102.1982032613868
$ws.Range("F12").Value = 8098.999999999988

$ws.Range("B13").Value = "P_1093"
$ws.Range("C13").Value = 2692.961134265342
$ws.Range("D13").Value = 123.3650092500005
$ws.Range("E13").Value = 76.84668256242304
$ws.Range("F13").Value = 5594.785000000023

$ws.Range("B14").Value = "P_1257"
$ws.Range("C14").Value = 2155.45450642888
$ws.Range("D14").Value = 88.96876707959395
$ws.Range("E14").Value = 53.72884433391435
$ws.Range("F14").Value = 4034.864720162991

$ws.Range("B15").Value = "P_1143"
$ws.Range("C15").Value = 2284.129538760402
$ws.Range("D15").Value = 100.4333399999999
$ws.Range("E15").Value = 59.91512060420938
$ws.Range("F15").Value = 4554.799999999996

$ws.Range("B16").Value = "P_1272"
$ws.Range("C16").Value = 3800.168724681668
$ws.Range("D16").Value = 156.9272040000004
$ws.Range("E16").Value = 94.16916893510829
$ws.Range("F16").Value = 7116.880000000015

$ws.Range("B17").Value = "P_1127"
$ws.Range("C17").Value = 4883.946985743577
$ws.Range("D17").Value = 204.9613649999997
$ws.Range("E17").Value = 123.2169825691001
$ws.Range("F17").Value = 9295.299999999983

$ws.Range("B18").Value = "P_1260"
$ws.Range("C18").Value = 3205.965502316036
$ws.Range("D18").Value = 145.7358503319236
$ws.Range("E18").Value = 79.99428898205711
$ws.Range("F18").Value = 6609.335615960253

$ws.Range("B19").Value = "P_1123"
$ws.Range("C19").Value = 3565.449696020921
$ws.Range("D19").Value = 138.0354095644854
$ws.Range("E19").Value = 74.6646107620076
$ws.Range("F19").Value = 6260.109277300924

$ws.Range("B20").Value = "P_1295"
$ws.Range("C20").Value = 1952.784491228811
$ws.Range("D20").Value = 88.51399199999986
$ws.Range("E20").Value = 54.93293609434996
$ws.Range("F20").Value = 4014.239999999993

$ws.Range("B21").Value = "P_1375"
$ws.Range("C21").Value = 2797.217238919427
$ws.Range("D21").Value = 120.1460400000001
$ws.Range("E21").Value = 61.20184795059098
$ws.Range("F21").Value = 5448.800000000005

$ws.Range("B22").Value = "P_1376"
$ws.Range("C22").Value = 2216.318762186473
$ws.Range("D22").Value = 93.68710783880661
$ws.Range("E22").Value = 47.30750649483431
$ws.Range("F22").Value = 4248.848428063791

$ws.Range("B23").Value = "P_1414"
$ws.Range("C23").Value = 1160.898817845534
$ws.Range("D23").Value = 45.86764546122762
$ws.Range("E23").Value = 28.6187456631727
$ws.Range("F23").Value = 2080.165327039801

$ws.Range("B24").Value = "P_1131"
$ws.Range("C24").Value = 2596.949245184677
$ws.Range("D24").Value = 108.7814700000004
$ws.Range("E24").Value = 61.0705177994389
$ws.Range("F24").Value = 4933.400000000017

$ws.Range("B25").Value = "P_1112"
$ws.Range("C25").Value = 2416.567308364345
$ws.Range("D25").Value = 96.61075200000002
$ws.Range("E25").Value = 44.59591799591247
$ws.Range("F25").Value = 4381.440000000001

$ws.Range("B26").Value = "P_1132"
$ws.Range("C26").Value = 1840.499932054379
$ws.Range("D26").Value = 70.53331950000009
$ws.Range("E26").Value = 41.72447031467404
$ws.Range("F26").Value = 3198.790000000004

$ws.Range("B27").Value = "P_1253"
$ws.Range("C27").Value = 1306.370234871559
$ws.Range("D27").Value = 54.93316500000002
$ws.Range("E27").Value = 35.12202847676254
$ws.Range("F27").Value = 2491.300000000001

$ws.Range("B28").Value = "P_1135"
$ws.Range("C28").Value = 2674.917566227386
$ws.Range("D28").Value = 118.6907400000003
$ws.Range("E28").Value = 71.40030787805004
$ws.Range("F28").Value = 5382.800000000015

$ws.Range("B29").Value = "P_1109"
$ws.Range("C29").Value = 3650.604025650016
$ws.Range("D29").Value = 158.2335562500003
$ws.Range("E29").Value = 98.31126052625729
$ws.Range("F29").Value = 7176.125000000015

$ws.Range("B30").Value = "P_1353"
$ws.Range("C30").Value = 1759.668366209476
$ws.Range("D30").Value = 64.09599683416663
$ws.Range("E30").Value = 40.58505108367701
$ws.Range("F30").Value = 2906.847928987148

$ws.Range("B31").Value = "P_1424"
$ws.Range("C31").Value = 461.5235818385182
$ws.Range("D31").Value = 19.65511432815084
$ws.Range("E31").Value = 12.26175685510235
$ws.Range("F31").Value = 891.3884049048

$ws.Range("B32").Value = "P_1427"
$ws.Range("C32").Value = 2049.729710923711
$ws.Range("D32").Value = 98.43651703802752
$ws.Range("E32").Value = 58.59265312081038
$ws.Range("F32").Value = 4464.241135511452

$ws.Range("B33").Value = "P_1217"
$ws.Range("C33").Value = 714.5686680257891
$ws.Range("D33").Value = 25.474365
$ws.Range("E33").Value = 11.32326814729268
$ws.Range("F33").Value = 1155.3

$ws.Range("B34").Value = "P_1271"
$ws.Range("C34").Value = 4392.28873545451
$ws.Range("D34").Value = 190.32678
$ws.Range("E34").Value = 100.2043873143121
$ws.Range("F34").Value = 8631.599999999999

$ws.Range("B35").Value = "P_1368"
$ws.Range("C35").Value = 1858.106304233767
$ws.Range("D35").Value = 74.43648098050134
$ws.Range("E35").Value = 46.05181330529004
$ws.Range("F35").Value = 3375.804126099834

$ws.Range("B36").Value = "P_1094"
$ws.Range("C36").Value = 2893.98085224203
$ws.Range("D36").Value = 114.314256
$ws.Range("E36").Value = 66.33131540430304
$ws.Range("F36").Value = 5184.319999999997

$ws.Range("B37").Value = "P_1285"
$ws.Range("C37").Value = 1926.757652453498
$ws.Range("D37").Value = 86.91228000000004
$ws.Range("E37").Value = 50.52206982505907
$ws.Range("F37").Value = 3941.600000000001

$ws.Range("B38").Value = "P_1288"
$ws.Range("C38").Value = 891.7917779102276
$ws.Range("D38").Value = 38.28530474999997
$ws.Range("E38").Value = 21.60099836996058
$ws.Range("F38").Value = 1736.294999999999

$ws.Range("B39").Value = "P_1357"
$ws.Range("C39").Value = 2670.161434560074
$ws.Range("D39").Value = 102.2198741718037
$ws.Range("E39").Value = 65.1908425176393
$ws.Range("F39").Value = 4635.821957904929

$ws.Range("B40").Value = "P_1280"
$ws.Range("C40").Value = 1493.333206700058
$ws.Range("D40").Value = 66.38461199999999
$ws.Range("E40").Value = 39.61473453095603
$ws.Range("F40").Value = 3010.639999999999

$ws.Range("B41").Value = "P_1281"
$ws.Range("C41").Value = 2737.677560424062
$ws.Range("D41").Value = 117.436095
$ws.Range("E41").Value = 68.50094724358814
$ws.Range("F41").Value = 5325.900000000001

$ws.Range("B42").Value = "P_1296"
$ws.Range("C42").Value = 1889.984949153536
$ws.Range("D42").Value = 88.03021499999994
$ws.Range("E42").Value = 52.77573966803654
$ws.Range("F42").Value = 3992.299999999997

$ws.Range("B43").Value = "P_1367"
$ws.Range("C43").Value = 1962.81372280169
$ws.Range("D43").Value = 83.50570956409045
$ws.Range("E43").Value = 50.09945291411204
$ws.Range("F43").Value = 3787.107009709317

$ws.Range("B44").Value = "P_1125a"
$ws.Range("C44").Value = 928.4025525895778
$ws.Range("D44").Value = 41.83201878493561
$ws.Range("E44").Value = 20.95130077938576
$ws.Range("F44").Value = 1897.143709067375

$ws.Range("B45").Value = "P_1087"
$ws.Range("C45").Value = 4785.751954858837
$ws.Range("D45").Value = 181.5042442499993
$ws.Range("E45").Value = 108.871073795229
$ws.Range("F45").Value = 8231.484999999966

$ws.Range("B46").Value = "P_1422"
$ws.Range("C46").Value = 1055.913256616354
$ws.Range("D46").Value = 49.29635530629912
$ws.Range("E46").Value = 28.5739729539983
$ws.Range("F46").Value = 2235.66237216776

$ws.Range("B47").Value = "P_1137"
$ws.Range("C47").Value = 1869.313744595796
$ws.Range("D47").Value = 78.35881930843867
$ws.Range("E47").Value = 41.2342618612266
$ws.Range("F47").Value = 3553.687950496084

$ws.Range("B48").Value = "P_1294"
$ws.Range("C48").Value = 1097.271109219816
$ws.Range("D48").Value = 35.58278763351645
$ws.Range("E48").Value = 18.68265108277637
$ws.Range("F48").Value = 1613.731865465599

$ws.Range("B49").Value = "P_1091"
$ws.Range("C49").Value = 8041.827617681056
$ws.Range("D49").Value = 325.5459794999994
$ws.Range("E49").Value = 201.0975713303817
$ws.Range("F49").Value = 14763.98999999997

$ws.Range("B50").Value = "P_1125"
$ws.Range("C50").Value = 1021.013616827905
$ws.Range("D50").Value = 37.52948423153453
$ws.Range("E50").Value = 21.83733499788648
$ws.Range("F50").Value = 1702.017425466418

$ws.Range("B51").Value = "P_1098"
$ws.Range("C51").Value = 2600.282604821728
$ws.Range("D51").Value = 121.5743120129884
$ws.Range("E51").Value = 58.89065695678865
$ws.Range("F51").Value = 5513.574240951853
